# Applies crypto price/volume updates to match the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '60.792.34'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.57%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.384.10'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -2.14%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '571.69'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.60%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '141.80'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -4.61%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '3.383.62'; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -2.25%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  -0.13%  '; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -3.56%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -0.91%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.392'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +0.16%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '3.964.44'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -2.26%  '; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'E15'; Value = '  +0.88%  '; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -2.90%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '3.386.05'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -2.14%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '60.924.48'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  -1.55%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '6.28'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -0.96%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '14.16'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '8.97'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -5.26%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '388.62'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.83%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  -1.47%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '73.50'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +1.05%  '; ForceText = $false }
    @{ Cell = 'E25'; Value = '  +0.16%  '; ForceText = $false }
    @{ Cell = 'E26'; Value = '  -3.78%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '3.523.66'; ForceText = $false }
    @{ Cell = 'E27'; Value = '  -2.10%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  -1.92%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  -0.12%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '7.42'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -5.11%  '; ForceText = $false }
    @{ Cell = 'B31'; Value = 'InternetComputer(DFINITY)'; ForceText = $false }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false }
    @{ Cell = 'D31'; Value = '8.08'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -1.95%  '; ForceText = $false }
    @{ Cell = 'B32'; Value = 'Fetch.AI'; ForceText = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; ForceText = $false }
    @{ Cell = 'D32'; Value = '1.46'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -4.22%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '23.77'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.75%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  -1.78%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '3.414.36'; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -2.02%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '166.69'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.31%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '5.03'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -3.67%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -3.13%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '0.0779'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -1.45%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '26.83'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +2.48%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.783'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -1.72%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.999'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -0.22%  '; ForceText = $false }
    @{ Cell = 'E45'; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '41.73'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -1.49%  '; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -2.03%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.546.19'; ForceText = $false }
    @{ Cell = 'E48'; Value = '  -2.14%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  -4.14%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '6.82'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -2.15%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '22.95'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -1.42%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
